$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain plain text so numeric-looking
# strings like "24.418.30" are not coerced into numbers by COM.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '24.418.30'
$ws.Range('E2').Value = '  +10.16%  '

$ws.Range('D3').Value = '1.677.11'
$ws.Range('E3').Value = '  +5.54%  '

$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.22%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '305.72'
$ws.Range('E5').Value = '  +2.42%  '

$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = '0.9967'
$ws.Range('E6').Value = '  +0.32%  '

$ws.Range('E7').Value = '  +1.59%  '

$ws.Range('D8').Value = '0.3424'
$ws.Range('E8').Value = '  +2.63%  '

$ws.Range('D9').Value = '47.96'
$ws.Range('E9').Value = '  +17.34%  '

$ws.Range('D10').Value = '1.160'
$ws.Range('E10').Value = '  +4.07%  '

$ws.Range('D11').Value = '0.07214'
$ws.Range('E11').Value = '  +3.84%  '

$ws.Range('D12').Value = '0.9975'
$ws.Range('E12').Value = '  -0.50%  '

$ws.Range('D13').Value = '6.119'
$ws.Range('E13').Value = '  +5.19%  '

$ws.Range('D14').Value = '20.13'
$ws.Range('E14').Value = '  +3.93%  '

$ws.Range('D15').Value = '6.721'
$ws.Range('E15').Value = '  +2.96%  '

$ws.Range('D16').Value = '1.675.46'
$ws.Range('E16').Value = '  +5.67%  '

$ws.Range('D17').Value = '0.00001100'
$ws.Range('E17').Value = '  +3.79%  '

$ws.Range('D18').Value = '0.9963'
$ws.Range('E18').Value = '  +0.15%  '

$ws.Range('D19').Value = '0.06653'
$ws.Range('E19').Value = '  +1.08%  '

$ws.Range('D20').Value = '80.61'
$ws.Range('E20').Value = '  +6.11%  '

$ws.Range('D21').Value = '16.43'
$ws.Range('E21').Value = '  +3.78%  '

$ws.Range('D22').Value = '6.095'
$ws.Range('E22').Value = '  +3.06%  '

$ws.Range('D23').Value = '12.12'
$ws.Range('E23').Value = '  +4.18%  '

$ws.Range('D24').Value = '24.384.40'
$ws.Range('E24').Value = '  +10.12%  '

$ws.Range('D25').Value = '2.416'
$ws.Range('E25').Value = '  +1.83%  '

$ws.Range('D26').Value = '2.649'
$ws.Range('E26').Value = '  +5.81%  '

$ws.Range('D27').Value = '153.30'
$ws.Range('E27').Value = '  +3.37%  '

$ws.Range('D28').Value = '19.43'
$ws.Range('E28').Value = '  +1.30%  '

$ws.Range('D29').Value = '1.857.27'
$ws.Range('E29').Value = '  +5.81%  '

$ws.Range('D30').Value = '127.68'
$ws.Range('E30').Value = '  +4.70%  '

$ws.Range('D31').Value = '6.262'
$ws.Range('E31').Value = '  +6.27%  '

$ws.Range('D32').Value = '4.022'
$ws.Range('E32').Value = '  +1.34%  '

$ws.Range('D33').Value = '0.9730'
$ws.Range('E33').Value = '  +5.57%  '

$ws.Range('D34').Value = '0.08426'
$ws.Range('E34').Value = '  +3.59%  '

$ws.Range('D35').Value = '1.695'
$ws.Range('E35').Value = '  +5.42%  '

$ws.Range('D36').Value = '12.34'
$ws.Range('E36').Value = '  +5.67%  '

$ws.Range('D37').Value = '0.06378'
$ws.Range('E37').Value = '  +6.66%  '

$ws.Range('D38').Value = '5.304'
$ws.Range('E38').Value = '  +3.45%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.02307'
$ws.Range('E39').Value = '  +6.03%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '8.632'
$ws.Range('E40').Value = '  +3.50%  '

$ws.Range('D41').Value = '1.244'
$ws.Range('E41').Value = '  +0.90%  '

$ws.Range('D42').Value = '0.2090'
$ws.Range('E42').Value = '  +5.38%  '

$ws.Range('D43').Value = '0.6081'
$ws.Range('E43').Value = '  +5.20%  '

$ws.Range('D44').Value = '0.9968'
$ws.Range('E44').Value = '  +0.12%  '

$ws.Range('D45').Value = '3.758'
$ws.Range('E45').Value = '  -0.28%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.5871'
$ws.Range('E46').Value = '  +5.84%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '12.84'
$ws.Range('E47').Value = '  -0.16%  '

$ws.Range('D48').Value = '125.54'
$ws.Range('E48').Value = '  +0.21%  '

$ws.Range('D49').Value = '2.010'
$ws.Range('E49').Value = '  +3.73%  '

$ws.Range('D50').Value = '0.07161'
$ws.Range('E50').Value = '  +6.86%  '

$ws.Range('D51').Value = '75.69'
$ws.Range('E51').Value = '  +4.57%  '

# Restore the default (un-styled) cell style now that the text values are set,
# matching the original workbook which had no explicit style on these cells.
$ws.Range("D2:E51").Style = "Normal"
